$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.181.13'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.15%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.658.27'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5162'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.07%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2644'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06268'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07780'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.473'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.658.47'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.885.35'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5453'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.88%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8099'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.95%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.195.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.612'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.08'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.008'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.56%  '
$ws.Range("E24").Value = '  +0.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.69'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1222'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.262'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.16'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.87%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.443'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.86%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05942'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.62%  '
$ws.Range("E31").Value = '  -1.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.550'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.262'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.585'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9610'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.426'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.773'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5669'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.069'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01590'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8529'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.004'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.010.68'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.49'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.799.87'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.82%  '
$ws.Range("E46").Value = '  +2.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.56'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.006'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.71%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.013'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05169'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4191'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.04%  '
